# Edit derived from the target commit:
#  1. The table on slide 6 gets a different built-in table style
#     ({FC336AA9-80BC-4ED3-80EF-7390DACB3414} -> {B89C0504-1EFE-4E75-B703-5DF655BB3084}).
#  2. The deck's theme colour scheme ("Integral") is replaced with the
#     stock "Office" colour palette (dk1/lt1 stay black/white; dk2, lt2,
#     accent1-6, hlink and folHlink all change).

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{B89C0504-1EFE-4E75-B703-5DF655BB3084}")
    }
}

# --- 2. Swap the theme colour scheme to the "Office" palette --------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
